# Correct incomplete sentence in the "Erläuterungen" sheet (B6):
# "...die eine Hälfte der Mietpreise kleiner als der Median..."
#   -> "...die eine Hälfte der Mietpreise ist kleiner als der Median..."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Erläuterungen")

$ws.Range("B6").Value = "Der Median ist der Wert, der die Mietpreise in zwei gleich grosse Hälften teilt, d.h. die eine Hälfte der Mietpreise ist kleiner als der Median, die andere Hälfte grösser."

# Matches the author's cursor ending up on B7 after editing B6.
$ws.Range("B7").Select() | Out-Null
